$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.794.26"
$ws.Range("E2").Value = "  +0.15%  "
$ws.Range("D3").Value = "3.805.32"
$ws.Range("E3").Value = "  +0.82%  "
$ws.Range("E4").Value = "  +0.24%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "604.26"
$ws.Range("E5").Value = "  +1.44%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "165.86"
$ws.Range("E6").Value = "  -0.63%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  +0.29%  "
$ws.Range("E10").Value = "  +1.00%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.31"
$ws.Range("E11").Value = "  +0.20%  "
$ws.Range("E12").Value = "  -0.74%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "35.98"
$ws.Range("E13").Value = "  +0.06%  "
$ws.Range("D14").Value = "4.446.70"
$ws.Range("E14").Value = "  +0.92%  "
$ws.Range("D15").Value = "3.792.59"
$ws.Range("E15").Value = "  +0.29%  "
$ws.Range("D16").Value = "67.820.21"
$ws.Range("E16").Value = "  +0.27%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.40"
$ws.Range("E17").Value = "  +0.14%  "
$ws.Range("E18").Value = "  +1.06%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "463.18"
$ws.Range("E20").Value = "  +1.34%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.86"
$ws.Range("E21").Value = "  -1.78%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.701"
$ws.Range("E22").Value = "  +1.09%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.0000146"
$ws.Range("E23").Value = "  -3.18%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "83.34"
$ws.Range("E24").Value = "  +0.18%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.09"
$ws.Range("E25").Value = "  +1.21%  "
$ws.Range("E26").Value = "  -0.53%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.01"
$ws.Range("E27").Value = "  -0.06%  "
$ws.Range("E28").Value = "  -0.09%  "
$ws.Range("D29").Value = "3.957.16"
$ws.Range("E29").Value = "  +0.95%  "
$ws.Range("E30").Value = "  +1.04%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.39"
$ws.Range("E31").Value = "  +2.45%  "
$ws.Range("E32").Value = "  -0.30%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "29.47"
$ws.Range("E33").Value = "  -0.60%  "
$ws.Range("E34").Value = "  +0.18%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.07"
$ws.Range("E35").Value = "  -0.37%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.1000"
$ws.Range("E36").Value = "  +0.12%  "
$ws.Range("E37").Value = "  +0.31%  "
$ws.Range("B38").Value = "Mantle"
$ws.Range("C38").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.996"
$ws.Range("E38").Value = "  +0.12%  "
$ws.Range("B39").Value = "Filecoin"
$ws.Range("C39").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.81"
$ws.Range("E39").Value = "  +1.13%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.20"
$ws.Range("E40").Value = "  -4.14%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("E41").Value = "  +0.09%  "
$ws.Range("E42").Value = "  -0.01%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "44.34"
$ws.Range("E43").Value = "  -1.70%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "47.71"
$ws.Range("E44").Value = "  -1.10%  "
$ws.Range("E45").Value = "  +0.62%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "27.87"
$ws.Range("E46").Value = "  +6.28%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "151.07"
$ws.Range("E47").Value = "  +1.06%  "
$ws.Range("E48").Value = "  +11.61%  "
$ws.Range("E49").Value = "  +0.62%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.85"
$ws.Range("E50").Value = "  +1.94%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "389.13"
$ws.Range("E51").Value = "  -0.19%  "
